# Append the latest scraped price-tracking row to the tracking sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Find the next empty row after the current data (row 37 is currently the last used row).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Write the new tracked values. A leading apostrophe forces Excel to store these
# as text (matching the Date/Price/Discount/Incredible columns, which are all
# stored as text in this sheet) instead of auto-converting them to a date/number.
$ws.Cells.Item($newRow, 1).Value = "'2026-02-07"
$ws.Cells.Item($newRow, 2).Value = "'530000"
$ws.Cells.Item($newRow, 3).Value = "'0"
$ws.Cells.Item($newRow, 4).Value = "'0"

# Drop the "quote prefix" text formatting that gets implicitly applied above so
# the new cells keep the same (default) style as the rest of the sheet.
$ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 4)).ClearFormats()
